$d = $word.ActiveDocument

# 1. "2021 & M.S." -> "2021, M.S."
$d.Content.Find.Execute("2021 & M.S.", $true, $false, $false, $false, $false, $true, 1, $false, "2021, M.S.", 2)

# 2. "Analyst on the J.P. Morgan Exotic Swaps Desk" -> "Worked as an analyst on the J.P. Morgan exotic rates desk"
$d.Content.Find.Execute("Analyst on the J.P. Morgan Exotic Swaps Desk", $true, $false, $false, $false, $false, $true, 1, $false, "Worked as an analyst on the J.P. Morgan exotic rates desk", 2)

# 3. "swap derivatives" -> "rates options"
$d.Content.Find.Execute("swap derivatives", $true, $false, $false, $false, $false, $true, 1, $false, "rates options", 2)

# 4. "sites" -> "destinations"
$d.Content.Find.Execute("sites", $true, $false, $false, $false, $false, $true, 1, $false, "destinations", 2)

# 5. " JS, Rust, C" -> " Rust, C"
$d.Content.Find.Execute(" JS, Rust, C", $true, $false, $false, $false, $false, $true, 1, $false, " Rust, C", 2)

# 6. "Poker, Home Improvement, Chess" -> "Home Improvement, Chess, Poker"
$d.Content.Find.Execute("Poker, Home Improvement, Chess", $true, $false, $false, $false, $false, $true, 1, $false, "Home Improvement, Chess, Poker", 2)
